$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "Metadata" sheet gains a new "Jurisdiction" property row, inserted
# right after "Contact" (row 11). Every row from the old "Description"
# row downward shifts down by one (old row 11 -> 12, 12 -> 13, 13 -> 14,
# 14 -> 15), and the sheet grows from 14 to 15 data rows.
#
# Capture the existing Property/Value pairs (rows 11-14) before they move,
# using .Text so long/special-character values (the Copyright blurb, the
# apostrophes in the Description, etc.) are carried over byte-for-byte
# without being retyped.
$a11 = $ws.Range("A11").Text
$b11 = $ws.Range("B11").Text
$a12 = $ws.Range("A12").Text
$b12 = $ws.Range("B12").Text
$a13 = $ws.Range("A13").Text
$b13 = $ws.Range("B13").Text
$a14 = $ws.Range("A14").Text
$b14 = $ws.Range("B14").Text

# New row 15 should look like the existing body rows (border/alignment),
# so clone the formatting from row 14 before anything is written into it.
$ws.Range("A14:B14").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)

# Push content down one row at a time, bottom-up, so nothing is clobbered.
$ws.Range("A15").Value = $a14
$ws.Range("B15").Value = $b14

$ws.Range("A14").Value = $a13
$ws.Range("B14").Value = $b13

$ws.Range("A13").Value = $a12
$ws.Range("B13").Value = $b12

$ws.Range("A12").Value = $a11
$ws.Range("B12").Value = $b11

# Row 11 becomes the new Jurisdiction property (value left blank).
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
